$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update existing weather-window values ---
$ws.Cells.Item(2, 2).Value = 45658          # B2 startDate
$ws.Cells.Item(2, 3).Value = 45658.25       # C2 endDate
$ws.Cells.Item(2, 4).Value = 2.6            # D2 minTemperature
$ws.Cells.Item(2, 5).Value = 3.4            # E2 maxTemperature
$ws.Cells.Item(2, 6).Value = 3.05           # F2 averageTemperature
$ws.Cells.Item(2, 7).Value = 3.1            # G2 medianTemperature

# --- Row 3: update existing weather-window values ---
$ws.Cells.Item(3, 2).Value = 45658.25       # B3 startDate
$ws.Cells.Item(3, 3).Value = 45658.5        # C3 endDate
$ws.Cells.Item(3, 4).Value = 2.8            # D3 minTemperature
$ws.Cells.Item(3, 5).Value = 3.9            # E3 maxTemperature
$ws.Cells.Item(3, 6).Value = 3.18           # F3 averageTemperature
$ws.Cells.Item(3, 7).Value = 3.2            # G3 medianTemperature

# --- Row 4: brand-new weather-window row ---
# Station id column mirrors the existing rows: a numeric-looking string kept
# as text (quote-prefix forces text storage just like the source rows).
$ws.Cells.Item(4, 1).Value = "'79049004"
$ws.Cells.Item(4, 1).NumberFormat = $ws.Cells.Item(3, 1).NumberFormat

# Match the date/number formatting already used by the data rows above.
$ws.Range("B4:C4").NumberFormat = $ws.Range("B3:C3").NumberFormat
$ws.Range("D4:G4").NumberFormat = $ws.Range("D3:G3").NumberFormat

$ws.Cells.Item(4, 2).Value = 45658.5        # B4 startDate
$ws.Cells.Item(4, 3).Value = 45658.75       # C4 endDate
$ws.Cells.Item(4, 4).Value = 3.5            # D4 minTemperature
$ws.Cells.Item(4, 5).Value = 4.4            # E4 maxTemperature
$ws.Cells.Item(4, 6).Value = 3.89           # F4 averageTemperature
$ws.Cells.Item(4, 7).Value = 3.8            # G4 medianTemperature

# Extend the "number stored as text" ignore hint (green-triangle suppression
# for the text-typed weatherStationId column) over the newly added row.
$ws.Range("A1:G4").Errors.Item(9).Ignore = $true
